$d = $word.ActiveDocument

# Locate the paragraph holding the old line "披金成王，伴坤远航".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*披金成王*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Grab the paragraph's own OOXML so we can reuse its exact <w:p ...>
    # attributes (paraId/textId/rsid/...) and its run's <w:rPr> without
    # hard-coding them.
    $openXml = $target.Range.WordOpenXML

    $pTag = "<w:p>"
    if ($openXml -match '<w:p( [^>]*)?>') {
        $pTag = $Matches[0]
    }

    $rPr = ""
    if ($openXml -match '<w:r>\s*(<w:rPr>.*?</w:rPr>)?\s*<w:t') {
        $rPr = $Matches[1]
    }

    # New text split across three runs: "披" / "龟" / "成王，伴坤远航"
    # (the middle character 金 becomes 龟).
    $newXml = $pTag +
              "<w:r>" + $rPr + "<w:t>披</w:t></w:r>" +
              "<w:r>" + $rPr + "<w:t>龟</w:t></w:r>" +
              "<w:r>" + $rPr + "<w:t>成王，伴坤远航</w:t></w:r>" +
              "</w:p>"

    # Replacing the whole paragraph range (including its end-of-paragraph
    # mark) drops the inherited paragraph-mark run properties (the old
    # <w:pPr><w:rPr>...</w:rPr></w:pPr>) while the three new runs keep the
    # original run formatting.
    $target.Range.InsertXML($newXml)
}
